$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddress, $value) {
    $rng = $ws.Range($cellAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '63.088.53'
$ws.Range("E2").Value = '  +2.07%  '

# Row 3
Set-TextValue "D3" '2.462.69'
$ws.Range("E3").Value = '  +2.03%  '

# Row 4
$ws.Range("E4").Value = '  -0.28%  '

# Row 5
Set-TextValue "D5" '576.83'
$ws.Range("E5").Value = '  +1.36%  '

# Row 6
Set-TextValue "D6" '146.45'
$ws.Range("E6").Value = '  +1.89%  '

# Row 7
$ws.Range("E7").Value = '  +0.17%  '

# Row 8
Set-TextValue "D8" '0.541'
$ws.Range("E8").Value = '  +1.03%  '

# Row 9
Set-TextValue "D9" '2.461.09'
$ws.Range("E9").Value = '  +1.37%  '

# Row 10
$ws.Range("E10").Value = '  +2.19%  '

# Row 11
$ws.Range("E11").Value = '  +1.67%  '

# Row 12
$ws.Range("E12").Value = '  +1.23%  '

# Row 13
$ws.Range("E13").Value = '  +2.27%  '

# Row 14
Set-TextValue "D14" '29.14'
$ws.Range("E14").Value = '  +9.74%  '

# Row 15
Set-TextValue "D15" '0.0000180'
$ws.Range("E15").Value = '  +2.75%  '

# Row 16
Set-TextValue "D16" '2.906.09'
$ws.Range("E16").Value = '  +2.52%  '

# Row 17
Set-TextValue "D17" '63.012.51'
$ws.Range("E17").Value = '  +2.21%  '

# Row 18
Set-TextValue "D18" '2.460.75'
$ws.Range("E18").Value = '  +1.49%  '

# Row 19
Set-TextValue "D19" '7.97'
$ws.Range("E19").Value = '  -0.27%  '

# Row 20
Set-TextValue "D20" '11.11'
$ws.Range("E20").Value = '  +3.81%  '

# Row 21
Set-TextValue "D21" '330.54'
$ws.Range("E21").Value = '  +1.83%  '

# Row 22
Set-TextValue "D22" '2.24'
$ws.Range("E22").Value = '  +8.92%  '

# Row 23
$ws.Range("E23").Value = '  +1.07%  '

# Row 24
$ws.Range("E24").Value = '  -0.02%  '

# Row 25
$ws.Range("E25").Value = '  +2.08%  '

# Row 26
Set-TextValue "D26" '663.93'
$ws.Range("E26").Value = '  +5.77%  '

# Row 27
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D27" '9.01'
$ws.Range("E27").Value = '  +6.88%  '

# Row 28
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue "D28" '1.11'
$ws.Range("E28").Value = '  +10.86%  '

# Row 29
Set-TextValue "D29" '0.0000101'
$ws.Range("E29").Value = '  +5.07%  '

# Row 30
$ws.Range("E30").Value = '  +2.39%  '

# Row 31
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D31" '1.45'
$ws.Range("E31").Value = '  +3.32%  '

# Row 32
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D32" '8.18'
$ws.Range("E32").Value = '  +1.22%  '

# Row 33
Set-TextValue "D33" '1.90'
$ws.Range("E33").Value = '  +4.77%  '

# Row 34
$ws.Range("E34").Value = '  +3.21%  '

# Row 35
Set-TextValue "D35" '1.55'
$ws.Range("E35").Value = '  +4.90%  '

# Row 36
$ws.Range("E36").Value = '  +0.11%  '

# Row 37
$ws.Range("E37").Value = '  +3.56%  '

# Row 38
Set-TextValue "D38" '5.53'
$ws.Range("E38").Value = '  +2.62%  '

# Row 39
Set-TextValue "D39" '153.42'
$ws.Range("E39").Value = '  +1.00%  '

# Row 40
Set-TextValue "D40" '0.374'
$ws.Range("E40").Value = '  +0.44%  '

# Row 41
Set-TextValue "D41" '18.88'
$ws.Range("E41").Value = '  +2.38%  '

# Row 42
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D42" '2.74'
$ws.Range("E42").Value = '  +6.23%  '

# Row 43
$ws.Range("B43").Value = 'BabyDogeCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue "D43" '0.0₆0344'
$ws.Range("E43").Value = '  +21.38%  '

# Row 44
$ws.Range("E44").Value = '  +3.31%  '

# Row 45
Set-TextValue "D45" '42.39'
$ws.Range("E45").Value = '  +0.85%  '

# Row 46
$ws.Range("E46").Value = '  +0.03%  '

# Row 47
Set-TextValue "D47" '15.12'
$ws.Range("E47").Value = '  +27.50%  '

# Row 48
Set-TextValue "D48" '146.75'
$ws.Range("E48").Value = '  +2.59%  '

# Row 49
Set-TextValue "D49" '3.64'
$ws.Range("E49").Value = '  +2.17%  '

# Row 50
$ws.Range("E50").Value = '  +3.84%  '

# Row 51
Set-TextValue "D51" '0.608'
